$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H6").Value = "['af', 'fr', 'en', 'ar']"
$ws.Range("I6").Value = "['en', 'tr', 'pl', 'de', 'ar', 'id', 'ca', 'fr', 'es', 'et', 'tl', 'undetected']"
$ws.Range("K6").Value = "['en', 'ja', 'pt', 'id', 'ar', 'fr']"
$ws.Range("R6").Value = "['en', 'nl', 'tr', 'so', 'id', 'ar', 'fr', 'es', 'undetected']"
$ws.Range("U6").Value = "['fa', 'en', 'so', 'pt', 'ca', 'ar', 'id', 'et', 'es', 'undetected']"
$ws.Range("W6").Value = "['undetected', 'en', 'de', 'ar']"
$ws.Range("Y6").Value = "['en', 'ca', 'ar', 'es', 'undetected']"
